$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-derive rows 2-4 (ECs x {ECs,FAPs,sCs}) and add new rows 5-10 so that the
# sending/target clusters form the full 3x3 cross of ECs/FAPs/sCs (per Dr Hou advice).
$rows = @(
    ,("ECs", "Angpt1", "Itga5", "ECs", 1, 0.3333333333333333, 0.07171233333333334, 0.215137, 0.004947717871829782, 0.004947717871829783, 3, 1, 28.72417333333333, 86.17251999999999, 0.4233259107972328, 0.4233259107972328, 2.059877492804444, 18.53889743524, 0.002094497174460089, 0.002094497174460089)
    ,("ECs", "Angpt1", "Itga5", "FAPs", 1, 0.3333333333333333, 0.07171233333333334, 0.215137, 0.004947717871829782, 0.004947717871829783, 3, 1, 30.56986233333333, 91.709587, 0.4505269713084062, 0.4505269713084062, 2.192236157602111, 19.730125418419, 0.002229080347683945, 0.002229080347683945)
    ,("ECs", "Angpt1", "Itga5", "sCs", 1, 0.3333333333333333, 0.07171233333333334, 0.215137, 0.004947717871829782, 0.004947717871829783, 3, 1, 8.559531999999999, 25.678596, 0.126147117894361, 0.126147117894361, 0.6138240119613333, 5.524416107652, 0.0006241403496857483, 0.0006241403496857484)
    ,("FAPs", "Angpt1", "Itga5", "ECs", 3, 1, 11.05178533333333, 33.155356, 0.7625064374239614, 0.7625064374239615, 3, 1, 28.72417333333333, 86.17251999999999, 0.4233259107972328, 0.4233259107972328, 317.4533975574577, 2857.080578017119, 0.3227887321112516, 0.3227887321112517)
    ,("FAPs", "Angpt1", "Itga5", "FAPs", 3, 1, 11.05178533333333, 33.155356, 0.7625064374239614, 0.7625064374239615, 3, 1, 30.56986233333333, 91.709587, 0.4505269713084062, 0.4505269713084062, 337.8515561775524, 3040.664005597972, 0.3435297158557801, 0.3435297158557802)
    ,("FAPs", "Angpt1", "Itga5", "sCs", 3, 1, 11.05178533333333, 33.155356, 0.7625064374239614, 0.7625064374239615, 3, 1, 8.559531999999999, 25.678596, 0.126147117894361, 0.126147117894361, 94.59811021779733, 851.3829919601759, 0.09618798945692963, 0.09618798945692965)
    ,("sCs", "Angpt1", "Itga5", "ECs", 3, 1, 3.370524666666666, 10.111574, 0.2325458447042087, 0.2325458447042088, 3, 1, 28.72417333333333, 86.17251999999999, 0.4233259107972328, 0.4233259107972328, 96.81553474960887, 871.3398127464799, 0.098442681511521, 0.09844268151152104)
    ,("sCs", "Angpt1", "Itga5", "FAPs", 3, 1, 3.370524666666666, 10.111574, 0.2325458447042087, 0.2325458447042088, 3, 1, 30.56986233333333, 91.709587, 0.4505269713084062, 0.4505269713084062, 103.0364750511042, 927.3282754599379, 0.1047681751049421, 0.1047681751049422)
    ,("sCs", "Angpt1", "Itga5", "sCs", 3, 1, 3.370524666666666, 10.111574, 0.2325458447042087, 0.2325458447042088, 3, 1, 8.559531999999999, 25.678596, 0.126147117894361, 0.126147117894361, 28.85011374112266, 259.6510236701039, 0.02933498808774558, 0.02933498808774558)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowVals = $rows[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
